# Auswertung Dozentenbefragung - add Thilo Skrotzki's survey response (row 19
# on Tabelle1) and move the active selection to reflect where the author was
# last working when they saved the file.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- Fill in the new respondent's row (row 19) on Tabelle1 ---------------
$ws1.Range("A19").Value = "Thilo Skrotzki"
$ws1.Range("B19").Value = "Professor"
$ws1.Range("C19").Value = "ja"
$ws1.Range("D19").Value = 10
$ws1.Range("E19").Value = 'ich musste etwas länger über Ihr Anliegen nachdenken und bin nach wie vor etwas unschlüssig bei der Antwort. Ich halte meine persönliche Erreichbarkeit für die Studierenden nicht für kompliziert. Zumindest habe ich nie derartige Rückmeldungen während meiner Zeit als Prüfungsausschussvorsitzender oder Studiendekan erhalten. Für mich sind auch Mailanfragen nicht störend, sondern eher erwünscht, weil ich dann die gesamte Terminkoordination viel besser im Überblick habe und keine Anfrage vergesse. Die möglichen Ansprechzeiten für ad hoc Gesprächswünsche sind m.M. nach auch einsehbar, VPIS gibt die Anwesenheit an der Hochschule an, natürlich im Rahmen der Lehrveranstaltungen. Die Antwort auf Ihre Frage wäre also "Ja" + 10 Min, weil ich sinngemäß schreiben würde " Tage laut VPIS Vorlesungsplan und nach Vereinbarung". Meine Zweifel sind: Wenn es Probleme mit der Erreichbarkeit von Lehrenden gibt, dann wird meiner Erfahrung nach solch ein Tool das nur wenig ändern. Ich möchte aber Ihre Verbesserungsansätze nicht im Keim ersticken.'

# Match the wrap-text style already used by the "Anmerkungen" column (E) on
# the neighbouring rows.
$ws1.Range("E19").WrapText = $true
$ws1.Rows.Item(19).RowHeight = 86.4

# --- Recalculate so the COUNTIF helper table / chart caches catch up -----
$excel.Calculate()

# --- Restore the view state recorded the last time the workbook was saved:
# Tabelle1 active, scrolled/selected at E19 -------------------------------
$ws1.Activate()
$ws1.Range("E19").Select()
